$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Main calibration data table (A2:B8), new row collected, old last row (9) removed ---
$ws.Range("A2").Value = 0.8309
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 2.0821
$ws.Range("B3").Value = 200

$ws.Range("A4").Value = 2.7664
$ws.Range("B4").Value = 300

$ws.Range("A5").Value = 3.4311
$ws.Range("B5").Value = 400

$ws.Range("A6").Value = 4.0469
$ws.Range("B6").Value = 500

$ws.Range("A7").Value = 4.6774
$ws.Range("B7").Value = 600

$ws.Range("A8").Value = 4.8729
$ws.Range("B8").Value = 620

# Old row 9 (previously the last data row) no longer exists - clear it
$ws.Range("A9").ClearContents()
$ws.Range("B9").ClearContents()

# --- Min/Max summary columns (M2:N3) mirror the new min/max of A/B ---
$ws.Range("M2").Value = 0.8309
$ws.Range("N2").Value = 0
$ws.Range("M3").Value = 4.8729
$ws.Range("N3").Value = 620

# --- Sample/voltage/pressure block for columns B:C moved up one row ---
# Row 22 now holds the "Sample" header (was row 23)
$ws.Range("B22").Value = "Sample"

# Row 23 now holds the voltage/pressure sub-headers (was row 24)
$ws.Range("B23").Value = "voltage"
$ws.Range("C23").Value = "pressure"

# Row 24 now holds the first voltage/pressure sample pair with updated coefficients
$ws.Range("B24").Value = 0.831
$ws.Range("C24").Formula = "=154.1*B24-125.6"

# Row 25 now holds the second voltage/pressure sample pair with updated coefficients
$ws.Range("B25").Value = 4.824
$ws.Range("C25").Formula = "=154.1*B25-125.6"

# Old row 26 B/C values are gone (table shrank by one row)
$ws.Range("B26").ClearContents()
$ws.Range("C26").ClearContents()

# --- N/O columns (rows 23-26) keep their row positions, only coefficients change ---
$ws.Range("O25").Formula = "=153.39*N25-127.45"
$ws.Range("O26").Formula = "=153.39*N26-127.45"

# --- Selection moved ---
$ws.Range("J26").Select() | Out-Null
